# Auto-generated Excel COM-interop script applying the scheduled market-data
# refresh described by the Phantom_Profits diff. Each block updates the price
# / profit columns (H..N) for a single leve row on one job sheet; the sheet
# order mirrors the workbook tab order (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 327.5
$ws.Range("I2").Value = 327.5
$ws.Range("K2").Value = 327.5
$ws.Range("M2").Value = -214.5

$ws.Range("H15").Value = 867.09753
$ws.Range("I15").Value = 867.09753
$ws.Range("K15").Value = 2601.29259
$ws.Range("M15").Value = -2432.29259

$ws.Range("H88").Value = 5051.1333
$ws.Range("J88").Value = 5326.7
$ws.Range("L88").Value = 5326.7
$ws.Range("N88").Value = -6138.7

$ws.Range("H91").Value = 5051.1333
$ws.Range("J91").Value = 5326.7
$ws.Range("L91").Value = 5326.7
$ws.Range("N91").Value = -8134.7

$ws.Range("H93").Value = 58724.375
$ws.Range("J93").Value = 58724.375
$ws.Range("L93").Value = 58724.375
$ws.Range("N93").Value = -63716.375

$ws.Range("H98").Value = 1933
$ws.Range("I98").Value = 1900
$ws.Range("J98").Value = 1949.5
$ws.Range("K98").Value = 1900
$ws.Range("L98").Value = 1949.5
$ws.Range("M98").Value = -402
$ws.Range("N98").Value = -4945.5

$ws.Range("H111").Value = 5127.091
$ws.Range("I111").Value = 5299.8887
$ws.Range("K111").Value = 15899.6661
$ws.Range("M111").Value = -12832.6661

$ws.Range("H122").Value = 1933
$ws.Range("I122").Value = 1900
$ws.Range("J122").Value = 1949.5
$ws.Range("K122").Value = 5700
$ws.Range("L122").Value = 5848.5
$ws.Range("M122").Value = -3250
$ws.Range("N122").Value = -10748.5

$ws.Range("H137").Value = 13890409
$ws.Range("I137").Value = 19608842
$ws.Range("K137").Value = 58826526
$ws.Range("M137").Value = -58823976

$ws.Range("H138").Value = 5100.3
$ws.Range("I138").Value = 7371.3335
$ws.Range("K138").Value = 22114.0005
$ws.Range("M138").Value = -16974.0005

$ws.Range("H141").Value = 2759.9
$ws.Range("I141").Value = 2371
$ws.Range("K141").Value = 7113
$ws.Range("M141").Value = -1933

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2174.9333
$ws.Range("I74").Value = 1749.7693
$ws.Range("J74").Value = 4938.5
$ws.Range("K74").Value = 1749.7693
$ws.Range("L74").Value = 4938.5
$ws.Range("M74").Value = -875.7692999999999
$ws.Range("N74").Value = -6686.5

$ws.Range("H77").Value = 2174.9333
$ws.Range("I77").Value = 1749.7693
$ws.Range("J77").Value = 4938.5
$ws.Range("K77").Value = 8748.8465
$ws.Range("L77").Value = 24692.5
$ws.Range("M77").Value = -4380.8465
$ws.Range("N77").Value = -33428.5

$ws.Range("H97").Value = 1873.875
$ws.Range("I97").Value = 1784.5714
$ws.Range("K97").Value = 1784.5714
$ws.Range("M97").Value = -1288.5714

$ws.Range("H110").Value = 1873.75
$ws.Range("I110").Value = 1833.3334
$ws.Range("K110").Value = 1833.3334
$ws.Range("M110").Value = 211.6666

$ws.Range("H117").Value = 59999.332

$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

$ws.Range("H132").Value = 4135.2
$ws.Range("I132").Value = 4203.5415
$ws.Range("K132").Value = 12610.6245
$ws.Range("M132").Value = -10080.6245

$ws.Range("H135").Value = 30852.666
$ws.Range("J135").Value = 30852.666
$ws.Range("L135").Value = 30852.666
$ws.Range("N135").Value = -40992.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 31525.4
$ws.Range("J2").Value = 31906.75
$ws.Range("L2").Value = 31906.75
$ws.Range("N2").Value = -32132.75

$ws.Range("H86").Value = 16705455
$ws.Range("I86").Value = 42096.69
$ws.Range("K86").Value = 42096.69
$ws.Range("M86").Value = -40973.69

$ws.Range("H89").Value = 16705455
$ws.Range("I89").Value = 42096.69
$ws.Range("K89").Value = 210483.45
$ws.Range("M89").Value = -204867.45

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H107").Value = 3926.4285
$ws.Range("I107").Value = 3926.4285
$ws.Range("K107").Value = 3926.4285
$ws.Range("M107").Value = -2006.4285

$ws.Range("H134").Value = 1555.6
$ws.Range("I134").Value = 1555.6
$ws.Range("K134").Value = 4666.799999999999
$ws.Range("M134").Value = -2131.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7672.5835
$ws.Range("I31").Value = 7370.091
$ws.Range("J31").Value = 11000
$ws.Range("K31").Value = 7370.091
$ws.Range("L31").Value = 11000
$ws.Range("M31").Value = -7075.091
$ws.Range("N31").Value = -11590

$ws.Range("H34").Value = 7672.5835
$ws.Range("I34").Value = 7370.091
$ws.Range("J34").Value = 11000
$ws.Range("K34").Value = 7370.091
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = -7168.091
$ws.Range("N34").Value = -11404

$ws.Range("H86").Value = 34998.5
$ws.Range("I86").Value = 34998
$ws.Range("J86").Value = 34999
$ws.Range("K86").Value = 34998
$ws.Range("L86").Value = 34999
$ws.Range("M86").Value = -33875
$ws.Range("N86").Value = -37245

$ws.Range("H89").Value = 34998.5
$ws.Range("I89").Value = 34998
$ws.Range("J89").Value = 34999
$ws.Range("K89").Value = 174990
$ws.Range("L89").Value = 174995
$ws.Range("M89").Value = -169374
$ws.Range("N89").Value = -186227

$ws.Range("H134").Value = 1965.5518
$ws.Range("I134").Value = 2013.2916
$ws.Range("J134").Value = 1736.4
$ws.Range("K134").Value = 6039.8748
$ws.Range("L134").Value = 5209.200000000001
$ws.Range("M134").Value = -3504.8748
$ws.Range("N134").Value = -10279.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2877799.5
$ws.Range("I4").Value = 857416.2
$ws.Range("K4").Value = 2572248.6
$ws.Range("M4").Value = -2572136.6

$ws.Range("H112").Value = 8424.75
$ws.Range("J112").Value = 9666.5
$ws.Range("L112").Value = 28999.5
$ws.Range("N112").Value = -31215.5

$ws.Range("H131").Value = 2375.3333
$ws.Range("I131").Value = 1799.8
$ws.Range("J131").Value = 3094.75
$ws.Range("K131").Value = 5399.4
$ws.Range("L131").Value = 9284.25
$ws.Range("M131").Value = -359.3999999999996
$ws.Range("N131").Value = -19364.25

$ws.Range("H138").Value = 12268.909
$ws.Range("I138").Value = 12268.909
$ws.Range("K138").Value = 36806.727
$ws.Range("M138").Value = -31666.727

$ws.Range("H139").Value = 3162.125
$ws.Range("I139").Value = 973.2
$ws.Range("K139").Value = 2919.6
$ws.Range("M139").Value = 2220.4

$ws.Range("H141").Value = 3796.2856
$ws.Range("I141").Value = 3796.2856
$ws.Range("K141").Value = 11388.8568
$ws.Range("M141").Value = -6208.856800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 50833.168
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 59399.8
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 59399.8
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -59711.8

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H100").Value = 43999.1
$ws.Range("J100").Value = 43999
$ws.Range("L100").Value = 43999
$ws.Range("N100").Value = -46163

$ws.Range("H113").Value = 6374.625
$ws.Range("I113").Value = 6330.6665
$ws.Range("J113").Value = 6506.5
$ws.Range("K113").Value = 6330.6665
$ws.Range("L113").Value = 6506.5
$ws.Range("M113").Value = -4160.6665
$ws.Range("N113").Value = -10846.5

$ws.Range("H126").Value = 8605.929
$ws.Range("I126").Value = 8348.6
$ws.Range("J126").Value = 9249.25
$ws.Range("K126").Value = 25045.8
$ws.Range("L126").Value = 27747.75
$ws.Range("M126").Value = -22575.8
$ws.Range("N126").Value = -32687.75

$ws.Range("H132").Value = 83335090
$ws.Range("I132").Value = 2333
$ws.Range("K132").Value = 6999
$ws.Range("M132").Value = -4469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 50004450
$ws.Range("I22").Value = 4749.6665
$ws.Range("J22").Value = 125004000
$ws.Range("K22").Value = 4749.6665
$ws.Range("L22").Value = 125004000
$ws.Range("M22").Value = -4454.6665
$ws.Range("N22").Value = -125004590

$ws.Range("H27").Value = 50004450
$ws.Range("I27").Value = 4749.6665
$ws.Range("J27").Value = 125004000
$ws.Range("K27").Value = 4749.6665
$ws.Range("L27").Value = 125004000
$ws.Range("M27").Value = -4642.6665
$ws.Range("N27").Value = -125004214

$ws.Range("H40").Value = 62504772
$ws.Range("I40").Value = 100003656
$ws.Range("J40").Value = 6633.1665
$ws.Range("K40").Value = 100003656
$ws.Range("L40").Value = 6633.1665
$ws.Range("M40").Value = -100003520
$ws.Range("N40").Value = -6905.1665

$ws.Range("H46").Value = 6799
$ws.Range("I46").Value = 6799
$ws.Range("K46").Value = 6799
$ws.Range("M46").Value = -6611

$ws.Range("H55").Value = 1030.3334
$ws.Range("I55").Value = 575.63635
$ws.Range("K55").Value = 575.63635
$ws.Range("M55").Value = -402.63635

$ws.Range("H60").Value = 22333
$ws.Range("I60").Value = 19999
$ws.Range("J60").Value = 23500
$ws.Range("K60").Value = 19999
$ws.Range("L60").Value = 23500
$ws.Range("M60").Value = -19490
$ws.Range("N60").Value = -24518

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 166701330
$ws.Range("I132").Value = 41601.8
$ws.Range("K132").Value = 124805.4
$ws.Range("M132").Value = -122275.4
